# Generate Report for Handoff
# - Drops the "handoff transform failed" row (the 1a45621d-...md file) from
#   every sheet, since that file no longer needs a row in the report.
# - Refreshes the remaining handoff-tracked file's name (now
#   a6e1d776-ca2c-4c1e-b047-befc8dde8234.md) plus its generated xlf
#   handoff-artifact name/hash and handoff timestamps on the zh-cn / de-de
#   sheets.

$wb = $excel.ActiveWorkbook

$oldMdName = "8da817aa-3de9-47c0-b83c-07f0d13c6835.md"
$newMdName = "a6e1d776-ca2c-4c1e-b047-befc8dde8234.md"

$oldZhXlf = "8da817aa-3de9-47c0-b83c-07f0d13c6835.b10c7a44e725e92666e2064d159bfb8395a03311.zh-cn.xlf"
$newZhXlf = "a6e1d776-ca2c-4c1e-b047-befc8dde8234.4fd030535bba4f8016cf207827d76bef0020e7fb.zh-cn.xlf"

$oldDeXlf = "8da817aa-3de9-47c0-b83c-07f0d13c6835.b10c7a44e725e92666e2064d159bfb8395a03311.de-de.xlf"
$newDeXlf = "a6e1d776-ca2c-4c1e-b047-befc8dde8234.4fd030535bba4f8016cf207827d76bef0020e7fb.de-de.xlf"

$newZhDatetime = "2016-02-17 04:41:31"
$newDeDatetime = "2016-02-17 04:41:41"

# ---------------------------------------------------------------------
# Overview sheet: drop row 3 (1a45621d-...md / "Handoff transform failed")
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(3).Delete()

$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6544b2fc4e733db531e8c539f9f3de0e304de0a7/e2e/$newMdName", [Type]::Missing, [Type]::Missing, $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6544b2fc4e733db531e8c539f9f3de0e304de0a7/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet: drop row 3, refresh the surviving row's file/hash/datetime
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows.Item(3).Delete()

$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhDatetime

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6544b2fc4e733db531e8c539f9f3de0e304de0a7/e2e/$newMdName", [Type]::Missing, [Type]::Missing, $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67a779bf8e5b753bcc890aaa6c780465858a60fd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newZhXlf", [Type]::Missing, [Type]::Missing, $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6544b2fc4e733db531e8c539f9f3de0e304de0a7/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet: drop row 3, refresh the surviving row's file/hash/datetime
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows.Item(3).Delete()

$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeDatetime

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6544b2fc4e733db531e8c539f9f3de0e304de0a7/e2e/$newMdName", [Type]::Missing, [Type]::Missing, $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/98f6b45be3c991473538a4e71a9c919624ed42ab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newDeXlf", [Type]::Missing, [Type]::Missing, $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6544b2fc4e733db531e8c539f9f3de0e304de0a7/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")
